$d = $word.ActiveDocument

# --- Step 1: swap the two English italic paragraphs (Objetivos <-> Programa resumido) ---
$null = $d.Content.Find.Execute("The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes", $true, $false, $false, $false, $false, $true, 1, $false, "§§TMP_SWAP_EN_BLOCK§§", 2)
$null = $d.Content.Find.Execute("To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues.", $true, $false, $false, $false, $false, $true, 1, $false, "The Administration of Organizations. 2 - The Administrative Process. 3 - Management Processes", 2)
$null = $d.Content.Find.Execute("§§TMP_SWAP_EN_BLOCK§§", $true, $false, $false, $false, $false, $true, 1, $false, "To present to the student the concept of an organization and the foundations of its administration; to characterize the various functional areas existing in the organizations; to awaken the interest of the students for management issues.", 2)

# --- Step 2: rewrite the "Docente(s) Responsavel(eis)" bullet paragraph (index 9) ---
$p9 = $d.Paragraphs.Item(9)
$rng9 = $d.Range($p9.Range.Start, $p9.Range.End)
$rng9.Text = "1 - A Administração das Organizações. 2 - O processo administrativo. 3 – Processos de Gestão" + [char]11 + "Apresentar ao aluno o conceito de uma organização e os fundamentos de sua administração;" + [char]11 + "Caracterizar as diversas áreas funcionais existentes nas organizações;" + [char]11 + "Despertar o interesse dos alunos para questões de gestão" + [char]11 + "1 - A Administração das organizações - definindo a administração" + [char]11 + "2 - O processo administrativo: planejamento, organização, direção, controle" + [char]11 + "3 – Processos de Gestão: Marketing, Finanças, Gestão de Pessoas, Produção e Operações, Pesquisa e Desenvolvimento, Tecnologia da Informação, Logística e Meio Ambiente." + [char]11 + "A disciplina será ministrada com duas estratégias pedagógicas a) aplicação de diferentes métodos ativos para compreender os principais conceitos necessários à gestão de negócios, e b) aplicação de conceitos por meio do Programa de Aprendizagem com Extensão, por meio do qual o alunos oferecem consultoria a micro e pequenas empresas da região de Lorena ou de parentes e amigos. Nestas consultorias times de alunos, orientados pelo professor, se debruçam sobre um pequeno problema de gestão da empresa e oferecem soluções." + [char]11 + "O sistema de avaliação será continuo com aplicação de avaliações escritas utilizando-se avaliações em papel como em sistema informacionais, ademais serão realizados seminários, projetos, entrega de trabalho em formato de artigo e Estudos de Casos." + [char]11 + "Avaliações em diversos formatos realizadas no decorrer do semestre. O peso maior da avaliação será aplicado ao Seminário Final da Disciplina, quando serão realizadas a apresentação oral do trabalho bem como a entrega do trabalho em formato de artigo; essa avaliação representará 70% da média do semestre." + [char]11 + "Os alunos em recuperação deverão realizar reuniões com o professor da disciplina, para orientar na execução de um trabalho em formato artigo científico em que se discutam as principais ferramentas de gestão e sua aplicação."

# --- Step 3: paragraph 11 becomes the bibliography references (moved from paragraph 14) ---
$p11 = $d.Paragraphs.Item(11)
$rng11 = $d.Range($p11.Range.Start, $p11.Range.End)
$rng11.Text = "LEMOS, Paulo de Mattos et al. Gestão estratégica de empresas. Rio de Janeiro: Fundação Getúlio Vargas, 2014." + [char]11 + "" + [char]11 + "Ludovico, Nelson. Gestão estratégica de negócios. São Paulo: Saraiva, 2018" + [char]11 + "" + [char]11 + "Serra, Fernando Ribeiro et al. Gestão estratégica: conceitos e casos. São Paulo: Atlas, 2014."

# --- Step 4: paragraph 14 becomes the "grupo social" text (moved from paragraph 17) ---
$p14 = $d.Paragraphs.Item(14)
$rng14 = $d.Range($p14.Range.Start, $p14.Range.End)
$rng14.Text = "O grupo social alvo da atividade é a comunidade com foco em emprendedores. O grupo social que participará da atividade serão proprietários de pequenos negócios que demandam conhecimentos de gestão."

# --- Step 5: rotate the three "Avaliacao" content runs (keep bold labels intact) ---
$null = $d.Content.Find.Execute("O grupo social alvo da atividade é a comunidade com foco em emprendedores. O grupo social que participará da atividade serão proprietários de pequenos negócios que demandam conhecimentos de gestão.", $true, $false, $false, $false, $false, $true, 1, $false, "§§TMP_SWAP_EN_BLOCK§§", 2)
$null = $d.Content.Find.Execute("- Contribuir para a gestão de organizações de pequeno e médio porte visando melhoria de rendas de comundades;" + [char]11 + "- contribuir para capacitar gestores de organizações de pequeno e medio porte.", $true, $false, $false, $false, $false, $true, 1, $false, "Grupos de alunos vão oferecer consultoria para micro e pequenas empresas, mentorados pelo professor, analisando e propondo melhorias em processos das empresas participantes visando aumentar a lucratividade e sustentabilidade dos negócios.", 2)
$null = $d.Content.Find.Execute("§§TMP_SWAP_EN_BLOCK§§", $true, $false, $false, $false, $false, $true, 1, $false, "- Contribuir para a gestão de organizações de pequeno e médio porte visando melhoria de rendas de comundades;" + [char]11 + "- contribuir para capacitar gestores de organizações de pequeno e medio porte.", 2)

# --- merge the old "Estabelecimento..." list onto the end of the "Norma de recuperacao" run ---
$null = $d.Content.Find.Execute("Grupos de alunos vão oferecer consultoria para micro e pequenas empresas, mentorados pelo professor, analisando e propondo melhorias em processos das empresas participantes visando aumentar a lucratividade e sustentabilidade dos negócios.", $true, $false, $false, $false, $false, $true, 1, $false, "Grupos de alunos vão oferecer consultoria para micro e pequenas empresas, mentorados pelo professor, analisando e propondo melhorias em processos das empresas participantes visando aumentar a lucratividade e sustentabilidade dos negócios." + [char]11 + "- Estabelecimento da comunicação aberta entre estudantes, grupo social e professor;" + [char]11 + "- Acompanhamento pelo professor e grupo social da atividade a ser desenvolvida pelos alunos;" + [char]11 + "- Exposição de cada grupo, sobre a proposta, desenvolvimento e finalização do projeto;" + [char]11 + "- Realização de avaliação conjunta dos resultados alcançados durante a atividade, incluindo benefícios   obtidos, lições aprendidas e desafios enfrentados;" + [char]11 + "- Conduzir sessões de discussão para revisar os resultados e identificar oportunidades de aplicação  prática;" + [char]11 + "- Apresentação do projeto final desenvolvido para grupo social;" + [char]11 + "- Avaliação do projeto apresentado, pelo grupo social e professor.", 2)

# --- Step 6: paragraph 19 becomes "849935 - Humberto Felipe da Silva" (moved out of paragraph 9) ---
$p19 = $d.Paragraphs.Item(19)
$rng19 = $d.Range($p19.Range.Start, $p19.Range.End)
$rng19.Text = "849935 - Humberto Felipe da Silva"

Write-Output "edit complete"
